$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E11").Value = 185
$ws.Range("E12").Value = 40
$ws.Range("E20").Value = 36
$ws.Range("E21").Value = 4
$ws.Range("E22").Value = 46
$ws.Range("E23").Value = 8
$ws.Range("E24").Value = 118
$ws.Range("E25").Value = 15
$ws.Range("E27").Value = 124
$ws.Range("E28").Value = 30
$ws.Range("E30").Value = 35
$ws.Range("E31").Value = 4
$ws.Range("E33").Value = 33
$ws.Range("E34").Value = 7
$ws.Range("E36").Value = 132
$ws.Range("E37").Value = 23
$ws.Range("E39").Value = 142
$ws.Range("E40").Value = 38
$ws.Range("E41").Value = 30
$ws.Range("E42").Value = 7
$ws.Range("E47").Value = 75
$ws.Range("E48").Value = 10
$ws.Range("E50").Value = 87
$ws.Range("E51").Value = 15
$ws.Range("E59").Value = 479
$ws.Range("E60").Value = 76
$ws.Range("E62").Value = 614
$ws.Range("E63").Value = 98
$ws.Range("E65").Value = 652
$ws.Range("E66").Value = 120
$ws.Range("E74").Value = 487
$ws.Range("E75").Value = 143
$ws.Range("E77").Value = 712
$ws.Range("E78").Value = 230
$ws.Range("E80").Value = 888
$ws.Range("E81").Value = 231
$ws.Range("E83").Value = 862
$ws.Range("E84").Value = 260
$ws.Range("E86").Value = 998
$ws.Range("E87").Value = 294
$ws.Range("E89").Value = 860
$ws.Range("E90").Value = 278
$ws.Range("E92").Value = 1102
$ws.Range("E93").Value = 332
$ws.Range("E95").Value = 906
$ws.Range("E96").Value = 335
$ws.Range("E98").Value = 1160
$ws.Range("E99").Value = 377
$ws.Range("E101").Value = 858
$ws.Range("E102").Value = 293
$ws.Range("E104").Value = 1136
$ws.Range("E105").Value = 365
$ws.Range("E107").Value = 839
$ws.Range("E108").Value = 399
$ws.Range("E110").Value = 1115
$ws.Range("E111").Value = 443
$ws.Range("E113").Value = 509
$ws.Range("E114").Value = 203
$ws.Range("E116").Value = 654
$ws.Range("E117").Value = 243
$ws.Range("E119").Value = 969
$ws.Range("E120").Value = 424
$ws.Range("E122").Value = 1250
$ws.Range("E123").Value = 546
$ws.Range("E125").Value = 941
$ws.Range("E126").Value = 542
$ws.Range("E131").Value = 576
$ws.Range("E132").Value = 256
$ws.Range("E134").Value = 775
$ws.Range("E135").Value = 329
$ws.Range("E137").Value = 724
$ws.Range("E138").Value = 405
$ws.Range("E140").Value = 925
$ws.Range("E141").Value = 512
$ws.Range("E143").Value = 779
$ws.Range("E144").Value = 447
$ws.Range("E146").Value = 918
$ws.Range("E147").Value = 554
$ws.Range("E149").Value = 538
$ws.Range("E150").Value = 471
$ws.Range("E152").Value = 707
$ws.Range("E153").Value = 622
$ws.Range("E155").Value = 691
$ws.Range("E156").Value = 692
$ws.Range("E158").Value = 850
$ws.Range("E159").Value = 847
$ws.Range("E161").Value = 551
$ws.Range("E162").Value = 547
$ws.Range("E164").Value = 672
$ws.Range("E165").Value = 632
$ws.Range("E167").Value = 755
$ws.Range("E168").Value = 420
$ws.Range("E170").Value = 967
$ws.Range("E171").Value = 497
$ws.Range("E173").Value = 534
$ws.Range("E174").Value = 387
$ws.Range("E176").Value = 700
$ws.Range("E177").Value = 450
$ws.Range("E182").Value = 603
$ws.Range("E183").Value = 512
$ws.Range("E185").Value = 537
$ws.Range("E186").Value = 365
$ws.Range("E188").Value = 704
$ws.Range("E189").Value = 495
$ws.Range("E191").Value = 645
$ws.Range("E192").Value = 498
$ws.Range("E194").Value = 850
$ws.Range("E195").Value = 642
$ws.Range("E197").Value = 724
$ws.Range("E198").Value = 503
$ws.Range("E200").Value = 914
$ws.Range("E201").Value = 619
$ws.Range("E203").Value = 289
$ws.Range("E204").Value = 340
$ws.Range("E206").Value = 339
$ws.Range("E207").Value = 399
$ws.Range("E209").Value = 588
$ws.Range("E210").Value = 811
$ws.Range("E211").Value = 30
$ws.Range("E212").Value = 744
$ws.Range("E213").Value = 1028
$ws.Range("E214").Value = 51
